# Update of Slovenia Prva Liga base data (30-05-2024 23:16)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename half-time goal columns ---
$ws.Range("I1").Value = "HTHG"
$ws.Range("J1").Value = "HTAG"

# --- Row 2 and Row 3 describe two matches that got reordered; swap their full content (columns B..AD) ---
$ws.Range("B2").Value = 6814327
$ws.Range("C2").Value = "Slovenia Prva Liga"
$ws.Range("D2").Value = 45129.52083333334
$ws.Range("E2").Value = "NS Mura"
$ws.Range("F2").Value = "NK Domzale"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "A"
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 3.3
$ws.Range("N2").Value = 3.4
$ws.Range("O2").Value = 1.909
$ws.Range("P2").Value = 3.4
$ws.Range("Q2").Value = 3.75
$ws.Range("R2").Value = -0.5
$ws.Range("S2").Value = 1.95
$ws.Range("T2").Value = 1.85
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.9
$ws.Range("W2").Value = 1.9
$ws.Range("X2").Value = -1
$ws.Range("Y2").Value = -1
$ws.Range("Z2").Value = 2.75
$ws.Range("AA2").Value = -1
$ws.Range("AB2").Value = 0.8500000000000001
$ws.Range("AC2").Value = 0.8999999999999999
$ws.Range("AD2").Value = -1

$ws.Range("B3").Value = 6816473
$ws.Range("C3").Value = "Slovenia Prva Liga"
$ws.Range("D3").Value = 45129.52083333334
$ws.Range("E3").Value = "NK Bravo"
$ws.Range("F3").Value = "NK Rogaska"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "H"
$ws.Range("L3").Value = 1.8
$ws.Range("M3").Value = 3.5
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 2.05
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 3.75
$ws.Range("R3").Value = -0.25
$ws.Range("S3").Value = 1.75
$ws.Range("T3").Value = 2.05
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.95
$ws.Range("W3").Value = 1.85
$ws.Range("X3").Value = 1.05
$ws.Range("Y3").Value = -1
$ws.Range("Z3").Value = -1
$ws.Range("AA3").Value = 0.75
$ws.Range("AB3").Value = -1
$ws.Range("AC3").Value = -0.5
$ws.Range("AD3").Value = 0.425

# --- Row 175 and Row 176 describe two matches that got reordered; swap their full content (columns B..AD) ---
$ws.Range("B175").Value = 7133777
$ws.Range("C175").Value = "Slovenia Prva Liga"
$ws.Range("D175").Value = 45430.41666666666
$ws.Range("E175").Value = "NK Radomlje"
$ws.Range("F175").Value = "NK Celje"
$ws.Range("G175").Value = 1
$ws.Range("H175").Value = 1
$ws.Range("I175").Value = 1
$ws.Range("J175").Value = 0
$ws.Range("K175").Value = "D"
$ws.Range("L175").Value = 3.05
$ws.Range("M175").Value = 3.5
$ws.Range("N175").Value = 2
$ws.Range("O175").Value = 2.9
$ws.Range("P175").Value = 3.6
$ws.Range("Q175").Value = 2.1
$ws.Range("R175").Value = 0.25
$ws.Range("S175").Value = 1.9
$ws.Range("T175").Value = 1.9
$ws.Range("U175").Value = 2.75
$ws.Range("V175").Value = 1.8
$ws.Range("W175").Value = 2
$ws.Range("X175").Value = -1
$ws.Range("Y175").Value = 2.6
$ws.Range("Z175").Value = -1
$ws.Range("AA175").Value = 0.45
$ws.Range("AB175").Value = -0.5
$ws.Range("AC175").Value = -1
$ws.Range("AD175").Value = 1

$ws.Range("B176").Value = 7124153
$ws.Range("C176").Value = "Slovenia Prva Liga"
$ws.Range("D176").Value = 45430.41666666666
$ws.Range("E176").Value = "NK Aluminij"
$ws.Range("F176").Value = "NK Domzale"
$ws.Range("G176").Value = 1
$ws.Range("H176").Value = 3
$ws.Range("I176").Value = 0
$ws.Range("J176").Value = 3
$ws.Range("K176").Value = "A"
$ws.Range("L176").Value = 2
$ws.Range("M176").Value = 3.6
$ws.Range("N176").Value = 3
$ws.Range("O176").Value = 1.333
$ws.Range("P176").Value = 4.75
$ws.Range("Q176").Value = 7
$ws.Range("R176").Value = -1.5
$ws.Range("S176").Value = 1.95
$ws.Range("T176").Value = 1.85
$ws.Range("U176").Value = 3.25
$ws.Range("V176").Value = 1.95
$ws.Range("W176").Value = 1.85
$ws.Range("X176").Value = -1
$ws.Range("Y176").Value = -1
$ws.Range("Z176").Value = 6
$ws.Range("AA176").Value = -1
$ws.Range("AB176").Value = 0.8500000000000001
$ws.Range("AC176").Value = 0.95
$ws.Range("AD176").Value = -1
